{"js": "// Replace the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table\n// with the new values from the commit. Cell text values are unique across\n// the whole document, so we walk the table's rows/cells in document order\n// and apply the Nth replacement to the Nth non-empty cell \u2014 this sidesteps\n// any issue where a *new* value happens to equal some *other* cell's *old*\n// value (e.g. replacement #11 below introduces \"63\u00f74=15, 3\", which is also\n// the ORIGINAL text of replacement #3 \u2014 a plain global text search/replace\n// run in a single pass could double-apply or mis-target).\nconst replacements = [\n  { before: \"35\u00f76=5, 5\", after: \"12\u00f74=3, 0\" },\n  { before: \"28\u00f76=4, 4\", after: \"45\u00f79=5, 0\" },\n  { before: \"63\u00f74=15, 3\", after: \"79\u00f78=9, 7\" },\n  { before: \"40\u00f72=20, 0\", after: \"61\u00f72=30, 1\" },\n  { before: \"26\u00f75=5, 1\", after: \"74\u00f77=10, 4\" },\n  { before: \"69\u00f73=23, 0\", after: \"52\u00f72=26, 0\" },\n  { before: \"20\u00f76=3, 2\", after: \"51\u00f76=8, 3\" },\n  { before: \"46\u00f73=15, 1\", after: \"40\u00f79=4, 4\" },\n  { before: \"19\u00f75=3, 4\", after: \"50\u00f72=25, 0\" },\n  { before: \"12\u00f78=1, 4\", after: \"10\u00f77=1, 3\" },\n  { before: \"34\u00f72=17, 0\", after: \"63\u00f74=15, 3\" },\n  { before: \"76\u00f75=15, 1\", after: \"69\u00f78=8, 5\" },\n  { before: \"50\u00f77=7, 1\", after: \"30\u00f79=3, 3\" },\n  { before: \"55\u00f74=13, 3\", after: \"53\u00f77=7, 4\" },\n  { before: \"17\u00f76=2, 5\", after: \"32\u00f76=5, 2\" },\n  { before: \"93\u00f72=46, 1\", after: \"11\u00f79=1, 2\" },\n  { before: \"99\u00f75=19, 4\", after: \"76\u00f72=38, 0\" },\n  { before: \"60\u00f73=20, 0\", after: \"59\u00f75=11, 4\" },\n  { before: \"49\u00f74=12, 1\", after: \"51\u00f75=10, 1\" },\n  { before: \"13\u00f79=1, 4\", after: \"88\u00f72=44, 0\" },\n  { before: \"82\u00f73=27, 1\", after: \"78\u00f76=13, 0\" },\n  { before: \"36\u00f75=7, 1\", after: \"63\u00f75=12, 3\" },\n  { before: \"53\u00f74=13, 1\", after: \"93\u00f74=23, 1\" },\n  { before: \"67\u00f73=22, 1\", after: \"80\u00f75=16, 0\" },\n  { before: \"32\u00f73=10, 2\", after: \"94\u00f78=11, 6\" },\n];\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet next = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (const cell of cells.items) {\n    if (next >= replacements.length) break;\n\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    await context.sync();\n\n    const text = para.text;\n    if (text === \"\") continue; // skip the blank spacer rows/cells\n\n    const expected = replacements[next];\n    if (text !== expected.before) {\n      throw new Error(\n        `Cell text mismatch at replacement #${next}: expected \"${expected.before}\", found \"${text}\"`\n      );\n    }\n\n    para.insertText(expected.after, \"Replace\");\n    await context.sync();\n    next++;\n  }\n}\n\nif (next !== replacements.length) {\n  throw new Error(`Only applied ${next} of ${replacements.length} replacements`);\n}\n", "ps1": "# Replace the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table\n# with the new values from the commit. Cell text values are unique across\n# the whole document, so we walk the table's rows/cells in document order\n# (Word COM is 1-based) and apply the Nth replacement to the Nth non-empty\n# cell \u2014 this sidesteps any issue where a *new* value happens to equal some\n# *other* cell's *old* value (e.g. replacement #11 below introduces\n# \"63\u00f74=15, 3\", which is also the ORIGINAL text of replacement #3 \u2014 a plain\n# global Find/Replace run in a single pass could double-apply or mis-target).\n$replacements = @(\n    @{ before = \"35\u00f76=5, 5\"; after = \"12\u00f74=3, 0\" },\n    @{ before = \"28\u00f76=4, 4\"; after = \"45\u00f79=5, 0\" },\n    @{ before = \"63\u00f74=15, 3\"; after = \"79\u00f78=9, 7\" },\n    @{ before = \"40\u00f72=20, 0\"; after = \"61\u00f72=30, 1\" },\n    @{ before = \"26\u00f75=5, 1\"; after = \"74\u00f77=10, 4\" },\n    @{ before = \"69\u00f73=23, 0\"; after = \"52\u00f72=26, 0\" },\n    @{ before = \"20\u00f76=3, 2\"; after = \"51\u00f76=8, 3\" },\n    @{ before = \"46\u00f73=15, 1\"; after = \"40\u00f79=4, 4\" },\n    @{ before = \"19\u00f75=3, 4\"; after = \"50\u00f72=25, 0\" },\n    @{ before = \"12\u00f78=1, 4\"; after = \"10\u00f77=1, 3\" },\n    @{ before = \"34\u00f72=17, 0\"; after = \"63\u00f74=15, 3\" },\n    @{ before = \"76\u00f75=15, 1\"; after = \"69\u00f78=8, 5\" },\n    @{ before = \"50\u00f77=7, 1\"; after = \"30\u00f79=3, 3\" },\n    @{ before = \"55\u00f74=13, 3\"; after = \"53\u00f77=7, 4\" },\n    @{ before = \"17\u00f76=2, 5\"; after = \"32\u00f76=5, 2\" },\n    @{ before = \"93\u00f72=46, 1\"; after = \"11\u00f79=1, 2\" },\n    @{ before = \"99\u00f75=19, 4\"; after = \"76\u00f72=38, 0\" },\n    @{ before = \"60\u00f73=20, 0\"; after = \"59\u00f75=11, 4\" },\n    @{ before = \"49\u00f74=12, 1\"; after = \"51\u00f75=10, 1\" },\n    @{ before = \"13\u00f79=1, 4\"; after = \"88\u00f72=44, 0\" },\n    @{ before = \"82\u00f73=27, 1\"; after = \"78\u00f76=13, 0\" },\n    @{ before = \"36\u00f75=7, 1\"; after = \"63\u00f75=12, 3\" },\n    @{ before = \"53\u00f74=13, 1\"; after = \"93\u00f74=23, 1\" },\n    @{ before = \"67\u00f73=22, 1\"; after = \"80\u00f75=16, 0\" },\n    @{ before = \"32\u00f73=10, 2\"; after = \"94\u00f78=11, 6\" },\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$next = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($next -ge $replacements.Count) {\n            continue\n        }\n\n        $cell = $t.Cell($r, $c)\n        $txt = $cell.Range.Text\n        # Drop the trailing cell mark (CR + BEL) Word always appends to\n        # Cell.Range.Text.\n        $txt = $txt.Substring(0, $txt.Length - 2)\n\n        if ($txt -eq \"\") {\n            continue\n        }\n\n        $expected = $replacements[$next]\n        if ($txt -ne $expected.before) {\n            throw \"Cell text mismatch at replacement $next (row $r, col $c): expected '$($expected.before)', found '$txt'\"\n        }\n\n        $cell.Range.Text = $expected.after\n        $next = $next + 1\n    }\n}\n\nif ($next -ne $replacements.Count) {\n    throw \"Only applied $next of $($replacements.Count) replacements\"\n}\n"}
